$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 3).Value = 45183
}
